$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "ТК_оригинал"
$ws.Range("O2").Value = "на момент выгрузки в элжуре"

for ($r = 4; $r -le 31; $r++) {
    $mVal = $ws.Cells.Item($r, 13).Value2
    $ws.Cells.Item($r, 15).Value = $mVal
    $ws.Cells.Item($r, 16).Formula = "=O$r-M$r"
}
